$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 16.506869
$ws.Range("H2").Value = 49.520607
$ws.Range("I2").Value = 0.2165594803671733
$ws.Range("J2").Value = 0.2165594803671733
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 20.02757633333333
$ws.Range("N2").Value = 60.082729
$ws.Range("O2").Value = 0.200201311135073
$ws.Range("P2").Value = 0.200201311135073
$ws.Range("Q2").Value = 330.5925789218336
$ws.Range("R2").Value = 2975.333210296503
$ws.Range("S2").Value = 0.0433554919082382
$ws.Range("T2").Value = 0.0433554919082382

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 16.506869
$ws.Range("H3").Value = 49.520607
$ws.Range("I3").Value = 0.2165594803671733
$ws.Range("J3").Value = 0.2165594803671733
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 33.29907966666667
$ws.Range("N3").Value = 99.89723900000001
$ws.Range("O3").Value = 0.3328670078646686
$ws.Range("P3").Value = 0.3328670078646686
$ws.Range("Q3").Value = 549.6635458782304
$ws.Range("R3").Value = 4946.971912904073
$ws.Range("S3").Value = 0.07208550625454842
$ws.Range("T3").Value = 0.07208550625454842

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.506869
$ws.Range("H4").Value = 49.520607
$ws.Range("I4").Value = 0.2165594803671733
$ws.Range("J4").Value = 0.2165594803671733
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 33.93321599999999
$ws.Range("N4").Value = 101.799648
$ws.Range("O4").Value = 0.3392060138062123
$ws.Range("P4").Value = 0.3392060138062122
$ws.Range("Q4").Value = 560.1311512607039
$ws.Range("R4").Value = 5041.180361346335
$ws.Range("S4").Value = 0.07345827808729355
$ws.Range("T4").Value = 0.07345827808729354

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.506869
$ws.Range("H5").Value = 49.520607
$ws.Range("I5").Value = 0.2165594803671733
$ws.Range("J5").Value = 0.2165594803671733
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 12.77731666666667
$ws.Range("N5").Value = 38.33195
$ws.Range("O5").Value = 0.1277256671940461
$ws.Range("P5").Value = 0.1277256671940461
$ws.Range("Q5").Value = 210.9134923881833
$ws.Range("R5").Value = 1898.22143149365
$ws.Range("S5").Value = 0.02766020411709313
$ws.Range("T5").Value = 0.02766020411709313

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 24.781512
$ws.Range("H6").Value = 74.34453600000001
$ws.Range("I6").Value = 0.3251174623990092
$ws.Range("J6").Value = 0.3251174623990092
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 20.02757633333333
$ws.Range("N6").Value = 60.082729
$ws.Range("O6").Value = 0.200201311135073
$ws.Range("P6").Value = 0.200201311135073
$ws.Range("Q6").Value = 496.3136232354161
$ws.Range("R6").Value = 4466.822609118744
$ws.Range("S6").Value = 0.06508894224518944
$ws.Range("T6").Value = 0.06508894224518944

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 24.781512
$ws.Range("H7").Value = 74.34453600000001
$ws.Range("I7").Value = 0.3251174623990092
$ws.Range("J7").Value = 0.3251174623990092
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 33.29907966666667
$ws.Range("N7").Value = 99.89723900000001
$ws.Range("O7").Value = 0.3328670078646686
$ws.Range("P7").Value = 0.3328670078646686
$ws.Range("Q7").Value = 825.2015423484562
$ws.Range("R7").Value = 7426.813881136106
$ws.Range("S7").Value = 0.1082208769133121
$ws.Range("T7").Value = 0.1082208769133121

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 24.781512
$ws.Range("H8").Value = 74.34453600000001
$ws.Range("I8").Value = 0.3251174623990092
$ws.Range("J8").Value = 0.3251174623990092
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 33.93321599999999
$ws.Range("N8").Value = 101.799648
$ws.Range("O8").Value = 0.3392060138062123
$ws.Range("P8").Value = 0.3392060138062122
$ws.Range("Q8").Value = 840.916399502592
$ws.Range("R8").Value = 7568.247595523328
$ws.Range("S8").Value = 0.110281798439159
$ws.Range("T8").Value = 0.110281798439159

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 24.781512
$ws.Range("H9").Value = 74.34453600000001
$ws.Range("I9").Value = 0.3251174623990092
$ws.Range("J9").Value = 0.3251174623990092
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 12.77731666666667
$ws.Range("N9").Value = 38.33195
$ws.Range("O9").Value = 0.1277256671940461
$ws.Range("P9").Value = 0.1277256671940461
$ws.Range("Q9").Value = 316.6412263028
$ws.Range("R9").Value = 2849.7710367252
$ws.Range("S9").Value = 0.04152584480134865
$ws.Range("T9").Value = 0.04152584480134863

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 24.67943933333333
$ws.Range("H10").Value = 74.038318
$ws.Range("I10").Value = 0.3237783348120013
$ws.Range("J10").Value = 0.3237783348120013
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.02757633333333
$ws.Range("N10").Value = 60.082729
$ws.Range("O10").Value = 0.200201311135073
$ws.Range("P10").Value = 0.200201311135073
$ws.Range("Q10").Value = 494.2693551122024
$ws.Range("R10").Value = 4448.424196009822
$ws.Range("S10").Value = 0.0648208471464933
$ws.Range("T10").Value = 0.0648208471464933

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 24.67943933333333
$ws.Range("H11").Value = 74.038318
$ws.Range("I11").Value = 0.3237783348120013
$ws.Range("J11").Value = 0.3237783348120013
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 33.29907966666667
$ws.Range("N11").Value = 99.89723900000001
$ws.Range("O11").Value = 0.3328670078646686
$ws.Range("P11").Value = 0.3328670078646686
$ws.Range("Q11").Value = 821.8026164893337
$ws.Range("R11").Value = 7396.223548404003
$ws.Range("S11").Value = 0.1077751255202757
$ws.Range("T11").Value = 0.1077751255202757

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 24.67943933333333
$ws.Range("H12").Value = 74.038318
$ws.Range("I12").Value = 0.3237783348120013
$ws.Range("J12").Value = 0.3237783348120013
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 33.93321599999999
$ws.Range("N12").Value = 101.799648
$ws.Range("O12").Value = 0.3392060138062123
$ws.Range("P12").Value = 0.3392060138062122
$ws.Range("Q12").Value = 837.4527456568959
$ws.Range("R12").Value = 7537.074710912064
$ws.Range("S12").Value = 0.1098275583083921
$ws.Range("T12").Value = 0.1098275583083921

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 24.67943933333333
$ws.Range("H13").Value = 74.038318
$ws.Range("I13").Value = 0.3237783348120013
$ws.Range("J13").Value = 0.3237783348120013
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 12.77731666666667
$ws.Range("N13").Value = 38.33195
$ws.Range("O13").Value = 0.1277256671940461
$ws.Range("P13").Value = 0.1277256671940461
$ws.Range("Q13").Value = 315.3370115177889
$ws.Range("R13").Value = 2838.0331036601
$ws.Range("S13").Value = 0.04135480383684011
$ws.Range("T13").Value = 0.0413548038368401

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 10.255437
$ws.Range("H14").Value = 30.766311
$ws.Range("I14").Value = 0.1345447224218162
$ws.Range("J14").Value = 0.1345447224218162
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 20.02757633333333
$ws.Range("N14").Value = 60.082729
$ws.Range("O14").Value = 0.200201311135073
$ws.Range("P14").Value = 0.200201311135073
$ws.Range("Q14").Value = 205.391547349191
$ws.Range("R14").Value = 1848.523926142719
$ws.Range("S14").Value = 0.02693602983515206
$ws.Range("T14").Value = 0.02693602983515206

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 10.255437
$ws.Range("H15").Value = 30.766311
$ws.Range("I15").Value = 0.1345447224218162
$ws.Range("J15").Value = 0.1345447224218162
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 33.29907966666667
$ws.Range("N15").Value = 99.89723900000001
$ws.Range("O15").Value = 0.3328670078646686
$ws.Range("P15").Value = 0.3328670078646686
$ws.Range("Q15").Value = 341.4966136794811
$ws.Range("R15").Value = 3073.46952311533
$ws.Range("S15").Value = 0.04478549917653235
$ws.Range("T15").Value = 0.04478549917653235

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 10.255437
$ws.Range("H16").Value = 30.766311
$ws.Range("I16").Value = 0.1345447224218162
$ws.Range("J16").Value = 0.1345447224218162
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 33.93321599999999
$ws.Range("N16").Value = 101.799648
$ws.Range("O16").Value = 0.3392060138062123
$ws.Range("P16").Value = 0.3392060138062122
$ws.Range("Q16").Value = 347.999958895392
$ws.Range("R16").Value = 3131.999630058528
$ws.Range("S16").Value = 0.04563837897136758
$ws.Range("T16").Value = 0.04563837897136759

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 10.255437
$ws.Range("H17").Value = 30.766311
$ws.Range("I17").Value = 0.1345447224218162
$ws.Range("J17").Value = 0.1345447224218162
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 12.77731666666667
$ws.Range("N17").Value = 38.33195
$ws.Range("O17").Value = 0.1277256671940461
$ws.Range("P17").Value = 0.1277256671940461
$ws.Range("Q17").Value = 131.03696610405
$ws.Range("R17").Value = 1179.33269493645
$ws.Range("S17").Value = 0.01718481443876421
$ws.Range("T17").Value = 0.01718481443876421
